$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (in-place run edits via Characters) ---
$ws.Range("A8").Characters(21, 2).Text = "47"
$ws.Range("C9").Characters(27, 10).Text = "11/20/2023"
$ws.Range("C9").Characters(48, 10).Text = "11/26/2023"

# --- Precinct crime-stat table updates (rows 14-30) ---
# Row 14
$ws.Range("F14").Value = 4
$ws.Range("G14").Value = 5
$ws.Range("H14").Value = -20
$ws.Range("I14").Value = 61
$ws.Range("K14").Value = -15.277777777777
$ws.Range("L14").Value = -28.235294117647
$ws.Range("M14").Value = -51.2
$ws.Range("N14").Value = -85.747663551401
# Row 15
$ws.Range("C15").Value = 4
$ws.Range("D15").Value = "'0"
$ws.Range("E15").Value = "***.*"
$ws.Range("G15").Value = 15
$ws.Range("H15").Value = -40
$ws.Range("I15").Value = 196
$ws.Range("K15").Value = -15.51724137931
$ws.Range("L15").Value = -2.970297029702
$ws.Range("M15").Value = -8.837209302325
$ws.Range("N15").Value = -64.684684684684
# Row 16
$ws.Range("C16").Value = 55
$ws.Range("D16").Value = 45
$ws.Range("E16").Value = 22.222222222222
$ws.Range("F16").Value = 225
$ws.Range("G16").Value = 200
$ws.Range("H16").Value = 12.5
$ws.Range("I16").Value = 2291
$ws.Range("J16").Value = 2339
$ws.Range("K16").Value = -2.052159042325
$ws.Range("L16").Value = 21.474019088017
$ws.Range("M16").Value = -31.242496998799
$ws.Range("N16").Value = -84.949415319931
# Row 17
$ws.Range("C17").Value = 84
$ws.Range("D17").Value = 73
$ws.Range("E17").Value = 15.068493150684
$ws.Range("F17").Value = 313
$ws.Range("G17").Value = 318
$ws.Range("H17").Value = -1.572327044025
$ws.Range("I17").Value = 3852
$ws.Range("J17").Value = 3800
$ws.Range("K17").Value = 1.368421052631
$ws.Range("L17").Value = 17.654245571166
$ws.Range("M17").Value = 27.634194831013
$ws.Range("N17").Value = -50.48206710374
# Row 18
$ws.Range("C18").Value = 43
$ws.Range("D18").Value = 33
$ws.Range("E18").Value = 30.30303030303
$ws.Range("F18").Value = 127
$ws.Range("G18").Value = 165
$ws.Range("H18").Value = -23.030303030303
$ws.Range("I18").Value = 1834
$ws.Range("J18").Value = 2179
$ws.Range("K18").Value = -15.832950894905
$ws.Range("L18").Value = -3.06553911205
$ws.Range("M18").Value = -37.427499147048
$ws.Range("N18").Value = -83.271002462829
# Row 19
$ws.Range("C19").Value = 79
$ws.Range("D19").Value = 72
$ws.Range("E19").Value = 9.722222222222
$ws.Range("F19").Value = 408
$ws.Range("G19").Value = 464
$ws.Range("H19").Value = -12.068965517241
$ws.Range("I19").Value = 5218
$ws.Range("J19").Value = 5442
$ws.Range("K19").Value = -4.116133774347
$ws.Range("L19").Value = 21.603355861104
$ws.Range("M19").Value = 32.807330109442
$ws.Range("N19").Value = -17.0956466476
# Row 20
$ws.Range("C20").Value = 36
$ws.Range("D20").Value = 34
$ws.Range("E20").Value = 5.882352941176
$ws.Range("F20").Value = 157
$ws.Range("G20").Value = 152
$ws.Range("H20").Value = 3.28947368421
$ws.Range("I20").Value = 1687
$ws.Range("J20").Value = 1683
$ws.Range("K20").Value = 0.237670825906
$ws.Range("L20").Value = 20.758768790264
$ws.Range("M20").Value = 29.669485011529
$ws.Range("N20").Value = -80.58241252302
# Row 21
$ws.Range("C21").Value = 303
$ws.Range("D21").Value = 257
$ws.Range("E21").Value = 17.898832684824
$ws.Range("F21").Value = 1243
$ws.Range("G21").Value = 1319
$ws.Range("H21").Value = -5.761940864291
$ws.Range("I21").Value = 15139
$ws.Range("J21").Value = 15747
$ws.Range("K21").Value = -3.861052898964
$ws.Range("L21").Value = 16.212481768634
$ws.Range("M21").Value = 1.939263349269
$ws.Range("N21").Value = -69.678944100622
# Row 22
$ws.Range("C22").Value = 8
$ws.Range("E22").Value = 14.285714285714
$ws.Range("F22").Value = 28
$ws.Range("G22").Value = 31
$ws.Range("H22").Value = -9.677419354838
$ws.Range("I22").Value = 259
$ws.Range("J22").Value = 318
$ws.Range("K22").Value = -18.553459119496
$ws.Range("L22").Value = 1.968503937007
$ws.Range("M22").Value = -34.595959595959
# Row 23
$ws.Range("C23").Value = 35
$ws.Range("D23").Value = 24
$ws.Range("E23").Value = 45.833333333333
$ws.Range("F23").Value = 100
$ws.Range("G23").Value = 109
$ws.Range("H23").Value = -8.256880733944
$ws.Range("I23").Value = 1415
$ws.Range("J23").Value = 1381
$ws.Range("K23").Value = 2.461984069514
$ws.Range("L23").Value = 5.75485799701
$ws.Range("M23").Value = 31.505576208178
# Row 24
$ws.Range("C24").Value = 189
$ws.Range("D24").Value = 220
$ws.Range("E24").Value = -14.090909090909
$ws.Range("F24").Value = 800
$ws.Range("G24").Value = 1000
$ws.Range("H24").Value = -20
$ws.Range("I24").Value = 11071
$ws.Range("J24").Value = 12195
$ws.Range("K24").Value = -9.216892168921
$ws.Range("L24").Value = 16.426543274792
$ws.Range("M24").Value = 16.15780086035
# Row 25
$ws.Range("C25").Value = 117
$ws.Range("D25").Value = 105
$ws.Range("E25").Value = 11.428571428571
$ws.Range("F25").Value = 472
$ws.Range("G25").Value = 479
$ws.Range("H25").Value = -1.461377870563
$ws.Range("I25").Value = 5582
$ws.Range("J25").Value = 5391
$ws.Range("K25").Value = 3.54294194027
$ws.Range("L25").Value = 28.056893782977
$ws.Range("M25").Value = -22.504512008885
# Row 26
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = 66.666666666666
$ws.Range("F26").Value = 19
$ws.Range("G26").Value = 25
$ws.Range("H26").Value = -24
$ws.Range("I26").Value = 305
$ws.Range("J26").Value = 347
$ws.Range("K26").Value = -12.103746397694
$ws.Range("L26").Value = -12.103746397694
# Row 27
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 6
$ws.Range("E27").Value = -66.666666666666
$ws.Range("F27").Value = 37
$ws.Range("G27").Value = 40
$ws.Range("H27").Value = -7.5
$ws.Range("I27").Value = 567
$ws.Range("J27").Value = 558
$ws.Range("K27").Value = 1.612903225806
$ws.Range("L27").Value = -9.28
# Row 28
$ws.Range("C28").Value = 9
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = 200
$ws.Range("F28").Value = 23
$ws.Range("H28").Value = 4.545454545454
$ws.Range("I28").Value = 223
$ws.Range("J28").Value = 317
$ws.Range("K28").Value = -29.652996845425
$ws.Range("L28").Value = -41.77545691906
$ws.Range("M28").Value = -52.350427350427
$ws.Range("N28").Value = -86.989498249708
# Row 29
$ws.Range("C29").Value = 7
$ws.Range("D29").Value = 3
$ws.Range("E29").Value = 133.333333333333
$ws.Range("F29").Value = 19
$ws.Range("G29").Value = 19
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 189
$ws.Range("J29").Value = 266
$ws.Range("K29").Value = -28.947368421052
$ws.Range("L29").Value = -39.423076923076
$ws.Range("M29").Value = -50.393700787401
$ws.Range("N29").Value = -87.743190661478
# Row 30
$ws.Range("C30").Value = 6
$ws.Range("D30").Value = "'0"
$ws.Range("E30").Value = "***.*"
$ws.Range("F30").Value = 14
$ws.Range("G30").Value = 10
$ws.Range("H30").Value = 40
$ws.Range("I30").Value = 73
$ws.Range("K30").Value = -13.095238095238
$ws.Range("L30").Value = 30.357142857142
